# Adds the new "Substrate vs. active-site residue" summary table (rows 20-26,
# columns B:H) to Sheet1, as part of the extra calculations added for the
# annual review (design_round_2 models summary).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - residue positions
$ws.Range("B20").Value = "Substrate"
$ws.Range("C20").Value = "Y65"
$ws.Range("D20").Value = "L68"
$ws.Range("E20").Value = "F81"
$ws.Range("F20").Value = "F83"
$ws.Range("G20").Value = "F89"
$ws.Range("H20").Value = "Y112"

# 5AM
$ws.Range("B21").Value = "5AM"
$ws.Range("C21").Value = "A/S"
$ws.Range("E21").Value = "Y"
$ws.Range("F21").Value = "F/H"
$ws.Range("G21").Value = "Y"
$ws.Range("H21").Value = "I/W"

# 5CP
$ws.Range("B22").Value = "5CP"
$ws.Range("C22").Value = "A/G"
$ws.Range("E22").Value = "Y"
$ws.Range("F22").Value = "F/N"
$ws.Range("G22").Value = "N/Y"
$ws.Range("H22").Value = "T/V/W"

# 5DZ
$ws.Range("B23").Value = "5DZ"
$ws.Range("C23").Value = "A/G"
$ws.Range("E23").Value = "F/Y"
$ws.Range("F23").Value = "A/G/H"
$ws.Range("G23").Value = "Y"
$ws.Range("H23").Value = "W"

# C35AM
$ws.Range("B24").Value = "C35AM"
$ws.Range("C24").Value = "A/G"
$ws.Range("E24").Value = "F/L/Y"
$ws.Range("F24").Value = "F/G/H"
$ws.Range("G24").Value = "Y"
$ws.Range("H24").Value = "T/W"

# Cyclo
$ws.Range("B25").Value = "Cyclo"
$ws.Range("C25").Value = "A/G"
$ws.Range("D25").Value = "A/L"
$ws.Range("E25").Value = "A/L/Y"
$ws.Range("F25").Value = "A/F/G"
$ws.Range("G25").Value = "A/Y"
$ws.Range("H25").Value = "I/W"

# M5
$ws.Range("B26").Value = "M5"
$ws.Range("C26").Value = "A"
$ws.Range("G26").Value = "Y"
$ws.Range("H26").Value = "W"

# Leave the new table selected, matching where the author's cursor ended up
$ws.Range("B20:H26").Select() | Out-Null
